$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update status for "High level Design Specs" row (row 8, column G)
# Change value from "Need to do review..." to "Ready for printing"
# and apply the same formatting used by the other "Ready for printing" cells (e.g. G3)
$ws.Range("G3").Copy()
$ws.Range("G8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("G8").Value = "Ready for printing"
$excel.CutCopyMode = $false

# Row 8 no longer needs the extra height for the long text - reset to default/auto height
$ws.Rows.Item(8).AutoFit()

# Update the active selection to G14
$ws.Range("G14").Select()
